# Update FoCus-results worksheet with new DPO row and revised figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in row 3 (Qwen2-5B-DPO-AVG)
$ws.Range("B3").Value = "0.85 ± 0.09"

# Update existing values in row 4 (Qwen2-5B-DPO-LENGTH-PRIOR)
$ws.Range("B4").Value = "0.6 ± 0.17"
$ws.Range("C4").Value = "0.8 ± 0.39"
$ws.Range("D4").Value = "0.03 ± 0.69"
$ws.Range("E4").Value = "0.34 ± 0.67"
$ws.Range("F4").Value = "0.4 ± 0.22"
$ws.Range("G4").Value = "0.149 ± 0.00"

# Add new row 5 (Qwen2-5B-DPO)
$ws.Range("A5").Value = "Qwen2-5B-DPO"
$ws.Range("B5").Value = "0.86 ± 0.09"
$ws.Range("C5").Value = "0.39 ± 0.48"
$ws.Range("D5").Value = "-0.47 ± 0.76"
$ws.Range("E5").Value = "0.16 ± 0.48"
$ws.Range("F5").Value = "0.22 ± 0.28"
$ws.Range("G5").Value = "0.595 ± 0.00"
